$wb = $excel.ActiveWorkbook
$wsProp = $wb.Worksheets.Item("Property1")
$wsCd = $wb.Worksheets.Item("Record_Cooldown")

# --- Fix up formatting first ------------------------------------------------
# A couple of cells in row 6 / G3 need to end up with a different direct
# cell format (border/number-format) than they started with. Do this BEFORE
# writing the new values: B6 currently carries a Text (@) number format, and
# assigning a numeric value into a Text-formatted cell stores it as a string,
# not a number - so the format has to be switched first.
$wsProp.Range("G5").Copy() | Out-Null
$wsProp.Range("G3").PasteSpecial(-4122) | Out-Null

$wsProp.Range("G4").Copy() | Out-Null
$wsProp.Range("G6").PasteSpecial(-4122) | Out-Null

$wsProp.Range("C6").Copy() | Out-Null
$wsProp.Range("B6").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Values ------------------------------------------------------------------
# "Public" row: Guild/G3 toggled off
$wsProp.Range("G3").Value = 0

# "Cache" row: now fully enabled across the row (new guild-related properties)
$wsProp.Range("B6:G6").Value = 1

# --- Data validation ----------------------------------------------------------
# Extend the TRUE/FALSE list validation so it also covers the newly-enabled
# B6:E6 cells (G3:G6 already carried it).
$wsProp.Range("B6:E6").Validation.Add(3, 1, 1, "TRUE,FALSE")

# --- Selection / active sheet -------------------------------------------------
# Selection within Property1 moves from A10 to G3, and Property1 becomes the
# active (selected) sheet/tab instead of Record_Cooldown.
$wsProp.Range("G3").Select()
$wsProp.Activate()

Write-Host "edit applied"
